# Weekly update: insert a new block of 5 price rows (week of 2021-09-09,
# serial 44448) for "Vega Modelo de Temuco - Palta" right before the
# existing data block that starts at row 635, pushing the old rows down
# by 5 (635-654 -> 640-659).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new blank rows at 635, shifting everything from 635 down
# (mirrors the old row 635's formatting, including the date style on D).
$ws.Range("A635:A639").EntireRow.Insert()

# Common columns shared by every row in this market/product block.
$mercadoId = 10
$mercado   = "Vega Modelo de Temuco"
$region    = "La Araucanía"
$fecha     = 44448
$codreg    = 9
$tipo      = "Fruta"
$productoId   = 100106
$producto     = "Oleaginosos"
$categoriaId  = 100106002
$categoria    = "Palta"
$variedad     = "Hass"
$kgUnidad     = 1

$rows = @(
    @{ Row = 635; Calidad = "1a nueva(o)";           Volumen = 100; PMin = 3000; PMax = 3000; PProm = 3000; Unidad = "$/kilo (en bandeja de 18 kilos)"; Origen = "Provincia de Limarí"; PrecioKg = 3000 },
    @{ Row = 636; Calidad = "2a nueva(o)";           Volumen = 200; PMin = 2500; PMax = 2500; PProm = 2500; Unidad = "$/kilo (en bandeja de 18 kilos)"; Origen = "Provincia de Limarí"; PrecioKg = 2500 },
    @{ Row = 637; Calidad = "Especial nueva (o)";    Volumen = 80;  PMin = 3500; PMax = 3500; PProm = 3500; Unidad = "$/kilo (en bandeja de 18 kilos)"; Origen = "Provincia de Limarí"; PrecioKg = 3500 },
    @{ Row = 638; Calidad = "Primera";                Volumen = 200; PMin = 2500; PMax = 2500; PProm = 2500; Unidad = "$/kilo (en caja de 8 kilos )";      Origen = "Perú";                  PrecioKg = 2500 },
    @{ Row = 639; Calidad = "Segunda";                Volumen = 300; PMin = 1800; PMax = 1800; PProm = 1800; Unidad = "$/kilo (en caja de 8 kilos )";      Origen = "Perú";                  PrecioKg = 1800 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = $r.Origen
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
